$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: section header for the "invalid secret key" test case
$ws.Range("A8").Value = "validateCreateCustomerAPIWithInValidSecretKey"

# Row 9: column headers (reuses existing shared strings)
$ws.Range("A9").Value = "name"
$ws.Range("B9").Value = "address"
$ws.Range("C9").Value = "description"
$ws.Range("D9").Value = "email"

# Column A: names
$ws.Range("A10").Value = "Kumar Kishan1"
$ws.Range("A11").Value = "Asif Khan1"
$ws.Range("A12").Value = "Maharaj Saxena1"
$ws.Range("A13").Value = "Kumar rawat1"

# Column B: addresses (reuses existing shared strings)
$ws.Range("B10").Value = "#30th main road"
$ws.Range("B11").Value = "#31th main road"
$ws.Range("B12").Value = "#32th main road"
$ws.Range("B13").Value = "#33th main road"

# Column C: descriptions (reuses existing shared strings)
$ws.Range("C10").Value = "This is request for kumar customer creation"
$ws.Range("C11").Value = "This is request for asif customer creation"
$ws.Range("C12").Value = "This is request for maharaj customer creation"
$ws.Range("C13").Value = "This is request for rawat customer creation"

# Column D: emails
$ws.Range("D10").Value = "kkinvalid123@gmail.com"
$ws.Range("D11").Value = "akinvalid123@gmail.com"
$ws.Range("D12").Value = "mhrjinvalid@gmail.com"
$ws.Range("D13").Value = "kmrinvalid@gmail.com"

# Hyperlinks on the email column for the newly added rows
[void]$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:kkinvalid123@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:akinvalid123@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:mhrjinvalid@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:kmrinvalid@gmail.com")

# Column width adjustments to fit the new (wider) content
$ws.Columns.Item(1).ColumnWidth = 40.83
$ws.Columns.Item(4).ColumnWidth = 20.83

# Restore selection to match the saved state of the workbook
[void]$ws.Range("D18").Select()
